$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet: bump the "Date" property to the new publication date.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# ---------------------------------------------------------------------------
# 2) Elements sheet: the two mapping columns (AK = col 37, AL = col 38) were
#    swapped - "Mapping: Spécification métier vers l'extension ROR
#    AvailableTimeTypeOfTime" now lives in AK (used to be AL) and
#    "Mapping: RIM Mapping" now lives in AL (used to be AK). The per-row data
#    that used to sit under each header moves along with the header.
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$headerAK = "Mapping: RIM Mapping"
$headerAL = "Mapping: Spécification métier vers l'extension ROR AvailableTimeTypeOfTime"

$oldAK = @{}
$oldAL = @{}
for ($r = 1; $r -le 6; $r++) {
    $oldAK[$r] = $elements.Cells.Item($r, 37).Value()
    $oldAL[$r] = $elements.Cells.Item($r, 38).Value()
}

for ($r = 1; $r -le 6; $r++) {
    $newAKValue = $oldAL[$r]
    $newALValue = $oldAK[$r]

    if ($newAKValue -eq $null) { $newAKValue = "" }
    if ($newALValue -eq $null) { $newALValue = "" }

    $elements.Cells.Item($r, 37).Value = $newAKValue
    $elements.Cells.Item($r, 38).Value = $newALValue
}

# Column widths follow the swapped content: AK (now the long French mapping
# text) becomes wide, AL (now the short "Mapping: RIM Mapping" / n/a values)
# becomes narrow.
$elements.Columns.Item(37).ColumnWidth = 81.16666666666667
$elements.Columns.Item(38).ColumnWidth = 24.166666666666668
